$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "Paveiksleliu api taisymas kad skelbimas rastu savo pics"
$ws.Range("C19").Value = 3.5

$ws.Range("C23").Select()
